$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# Insert 6 new rows right after the current last data row (row 36),
# carrying down formatting the same way Excel does when you insert rows
# in the middle of a formatted table.
$ws.Range("A37:O42").Insert()

# --- Row 37 ---
$ws.Range("A37").Value = 43567
$ws.Range("B37").Value = 242.76
$ws.Range("C37").Value = "Gasolina"
$ws.Range("D37").Value = "Gasolina"
$ws.Range("E37").Value = "Gasto"
$ws.Range("F37").Value = "Tarjeta Banamex"
$ws.Range("G37").Value = "Costco"
$ws.Range("K37").Formula = "=K36-B37"
$ws.Range("L37").Value = 32.5
$ws.Range("M37").Value = 504
$ws.Range("N37").Formula = "=SUM(K37:M37)"
$ws.Range("O37").Formula = "=N37-4000"

# --- Row 38 ---
$ws.Range("A38").Value = 43568
$ws.Range("B38").Value = 80
$ws.Range("C38").Value = "Café"
$ws.Range("D38").Value = "Comida"
$ws.Range("E38").Value = "Gasto"
$ws.Range("F38").Value = "Tarjeta Banamex"
$ws.Range("G38").Value = "Cafetería"
$ws.Range("K38").Formula = "=K37-B38"
$ws.Range("L38").Value = 32.5
$ws.Range("M38").Value = 504
$ws.Range("N38").Formula = "=SUM(K38:M38)"
$ws.Range("O38").Formula = "=N38-4000"

# --- Row 39 ---
$ws.Range("A39").Value = 43568
$ws.Range("B39").Value = 12
$ws.Range("C39").Value = "Estacionamiento Liverpool"
$ws.Range("D39").Value = "Estacionamiento"
$ws.Range("E39").Value = "Gasto"
$ws.Range("F39").Value = "Efectivo"
$ws.Range("G39").Value = "Galerías Celaya"
$ws.Range("K39").Value = 7035.88
$ws.Range("L39").Value = 32.5
$ws.Range("M39").Formula = "=M38-B39"
$ws.Range("N39").Formula = "=SUM(K39:M39)"
$ws.Range("O39").Formula = "=N39-4000"

# --- Row 40 ---
$ws.Range("A40").Value = 43568
$ws.Range("B40").Value = 30
$ws.Range("C40").Value = "Estacionamiento Centro Celaya"
$ws.Range("D40").Value = "Estacionamiento"
$ws.Range("E40").Value = "Gasto"
$ws.Range("F40").Value = "Efectivo"
$ws.Range("G40").Value = "Celaya Centro"
$ws.Range("K40").Value = 7035.88
$ws.Range("L40").Value = 32.5
$ws.Range("M40").Formula = "=M39-B40"
$ws.Range("N40").Formula = "=SUM(K40:M40)"
$ws.Range("O40").Formula = "=N40-4000"

# --- Row 41 ---
$ws.Range("A41").Value = 43568
$ws.Range("B41").Value = 25
$ws.Range("C41").Value = "Frappé Mexicano"
$ws.Range("D41").Value = "Comida"
$ws.Range("E41").Value = "Gasto"
$ws.Range("F41").Value = "Efectivo"
$ws.Range("G41").Value = "Celaya Centro"
$ws.Range("K41").Value = 7035.88
$ws.Range("L41").Value = 32.5
$ws.Range("M41").Formula = "=M40-B41"
$ws.Range("N41").Formula = "=SUM(K41:M41)"
$ws.Range("O41").Formula = "=N41-4000"

# --- Row 42 ---
$ws.Range("A42").Value = 43568
$ws.Range("B42").Value = 5
$ws.Range("C42").Value = "Puerquito de Barro"
$ws.Range("D42").Value = "Misc"
$ws.Range("E42").Value = "Gasto"
$ws.Range("F42").Value = "Efectivo"
$ws.Range("G42").Value = "Celaya Centro"
$ws.Range("K42").Value = 7035.88
$ws.Range("L42").Value = 32.5
$ws.Range("M42").Formula = "=M41-B42"
$ws.Range("N42").Formula = "=SUM(K42:M42)"
$ws.Range("O42").Formula = "=N42-4000"

# Keep selection consistent with where Excel would leave the cursor after
# typing the last new row of data.
[void]$ws.Range("Q42").Select()
